$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update actual hours burned down for Day 2 (column D) for the tasks
# that had work logged against them, per the burndown update.
$ws.Range("D7").Value = 0.75
$ws.Range("D9").Value = 0.25
$ws.Range("D12").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("D19").Value = 2
$ws.Range("D21").Value = 1
$ws.Range("D22").Value = 1
